$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text  = "919×9="
$t.Cell(1,2).Range.Text  = "999×9="
$t.Cell(1,3).Range.Text  = "937×6="
$t.Cell(1,4).Range.Text  = "686×2="
$t.Cell(1,5).Range.Text  = "727×9="

$t.Cell(5,1).Range.Text  = "371×8="
$t.Cell(5,2).Range.Text  = "597×2="
$t.Cell(5,3).Range.Text  = "429×7="
$t.Cell(5,4).Range.Text  = "686×6="
$t.Cell(5,5).Range.Text  = "325×2="

$t.Cell(10,1).Range.Text = "417×7="
$t.Cell(10,2).Range.Text = "254×5="
$t.Cell(10,3).Range.Text = "803×9="
$t.Cell(10,4).Range.Text = "973×2="
$t.Cell(10,5).Range.Text = "441×6="

$t.Cell(15,1).Range.Text = "654×6="
$t.Cell(15,2).Range.Text = "102×5="
$t.Cell(15,3).Range.Text = "136×2="
$t.Cell(15,4).Range.Text = "245×2="
$t.Cell(15,5).Range.Text = "272×3="

$t.Cell(20,1).Range.Text = "753×2="
$t.Cell(20,2).Range.Text = "373×5="
$t.Cell(20,3).Range.Text = "649×7="
$t.Cell(20,4).Range.Text = "195×2="
$t.Cell(20,5).Range.Text = "935×2="
